$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (rId1 / index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 4945
$ws1.Range("F3").Value = 2771
$ws1.Range("F5").Value = 2838
$ws1.Range("F9").Value = 1748
$ws1.Range("F10").Value = 758
$ws1.Range("F11").Value = 499
$ws1.Range("F13").Value = 435
$ws1.Range("F14").Value = 1081
$ws1.Range("F16").Value = 10
$ws1.Range("F17").Value = 95
$ws1.Range("F18").Value = 80
$ws1.Range("F19").Value = 1048
$ws1.Range("F22").Value = 678
$ws1.Range("F23").Value = 758
$ws1.Range("F24").Value = 155
$ws1.Range("F25").Value = 13
$ws1.Range("F27").Value = 561
$ws1.Range("F28").Value = 59
$ws1.Range("F29").Value = 1666
$ws1.Range("F30").Value = 1672
$ws1.Range("F31").Value = 412
$ws1.Range("F33").Value = 1587
$ws1.Range("F34").Value = 226
$ws1.Range("F35").Value = 2410
$ws1.Range("F36").Value = 417
$ws1.Range("F38").Value = 629
$ws1.Range("F39").Value = 119
$ws1.Range("F40").Value = 71
$ws1.Range("F42").Value = 823
$ws1.Range("F43").Value = 1519
$ws1.Range("F44").Value = 235
$ws1.Range("F47").Value = 71
$ws1.Range("F49").Value = 121

# Sheet 2: 演出 (rId2 / index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 106
$ws2.Range("F5").Value = 9

# Sheet 4: 全部类型 (rId4 / index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 4945
$ws4.Range("F3").Value = 2771
$ws4.Range("F4").Value = 2838
$ws4.Range("F5").Value = 1748
$ws4.Range("F8").Value = 758
$ws4.Range("F9").Value = 499
$ws4.Range("F11").Value = 435
$ws4.Range("F12").Value = 1081
$ws4.Range("F14").Value = 95
$ws4.Range("F15").Value = 1048
$ws4.Range("F17").Value = 678
$ws4.Range("F18").Value = 758
$ws4.Range("F19").Value = 155
$ws4.Range("F20").Value = 106
$ws4.Range("F21").Value = 106
$ws4.Range("F22").Value = 9
$ws4.Range("F23").Value = 13
$ws4.Range("F26").Value = 561
$ws4.Range("F27").Value = 1666
$ws4.Range("F28").Value = 1672
$ws4.Range("F29").Value = 412
$ws4.Range("F33").Value = 2410
$ws4.Range("F34").Value = 417
$ws4.Range("F40").Value = 119
$ws4.Range("F41").Value = 71
$ws4.Range("F43").Value = 823
$ws4.Range("F44").Value = 1519
$ws4.Range("F46").Value = 235
$ws4.Range("F48").Value = 71

